# plantilla_ejemplo.xlsx: rename "download" folder references to "document"
# (implicit in the workbook move) and refresh the example data block from a
# single "Sede Chile" example row into four example rows (Colombia, Bélgica,
# Argentina, Venezuela), dropping the now-unused "Estado" column (E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Estado" column (E) is removed entirely from the template.
$ws.Columns("E").Delete()

# Header row: A1 used to hold "Sede Chile"'s label; it becomes "Nombre_sede".
$ws.Range("A1").Value = "Nombre_sede"

# Seed rows 3:5 with row 2's number formats/styles before filling them in,
# so the new example rows render the same as the existing one (wrapped text
# for B, date formatting for C/D).
$ws.Range("A2:D2").Copy() | Out-Null
$ws.Range("A3:D5").PasteSpecial(-4122) | Out-Null

# Column A - Nombre_sede examples.
$ws.Range("A2").Value = "Sede Colombia"
$ws.Range("A5").Value = "Sede Bélgica"
$ws.Range("A3").Value = "Sede Argentina"
$ws.Range("A4").Value = "Sede Venezuela"

# Column B - Descripcion examples.
$ws.Range("B2").Value = "Proceso Ejemplo sede U colombia"
$ws.Range("B3").Value = "Proceso Ejemplo sede U argentina"
$ws.Range("B4").Value = "Proceso Ejemplo sede U venezuela"
$ws.Range("B5").Value = "Proceso Ejemplo sede U bélgica"

# Column C - Fecha_inicio examples.
$ws.Range("C2").Value = 45108
$ws.Range("C3").Value = 45109
$ws.Range("C4").Value = 45110
$ws.Range("C5").Value = 45111

# Column D - Fecha_fin examples (D2 keeps its original value, 45117).
$ws.Range("D3").Value = 45122
$ws.Range("D4").Value = 45127
$ws.Range("D5").Value = 45132

# Match the author's last selection before saving.
$ws.Range("E8").Select() | Out-Null
